# Update the CDA Logical Model "Metadata" sheet for ST.r2b:
#  - bump the Version and Date values
#  - insert a new "Jurisdiction" property row (empty value) right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version (row 3) and Date (row 8) values in place.
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row 11 for "Jurisdiction" (pushing Description, Purpose, ... down by one row),
# matching the existing row formatting.
$ws.Rows.Item(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
